$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend header row (row 1) with two more columns: P1=14, Q1=15
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Copy the header formatting (bold, border, centered) from O1 onto P1:Q1
$ws.Range("O1").Copy() | Out-Null
$ws.Range("P1:Q1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# For data rows 2-25: swap I<->K, swap M<->O, and add new P/Q columns = 2
for ($r = 2; $r -le 25; $r++) {
    $iVal = $ws.Cells.Item($r, 9).Value2   # column I
    $kVal = $ws.Cells.Item($r, 11).Value2  # column K
    $ws.Cells.Item($r, 9).Value = $kVal
    $ws.Cells.Item($r, 11).Value = $iVal

    $mVal = $ws.Cells.Item($r, 13).Value2  # column M
    $oVal = $ws.Cells.Item($r, 15).Value2  # column O
    $ws.Cells.Item($r, 13).Value = $oVal
    $ws.Cells.Item($r, 15).Value = $mVal

    $ws.Cells.Item($r, 16).Value = 2      # column P
    $ws.Cells.Item($r, 17).Value = 2      # column Q
}
